$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at Excel row 227, shifting existing rows 227-332 down to 228-333.
$ws.Rows("227").Insert()

# Populate the newly inserted row 227 with the new market record.
$ws.Cells.Item(227, 1).Value = 5
$ws.Cells.Item(227, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(227, 3).Value = 'Maule'

$newDate = Get-Date -Year 2021 -Month 11 -Day 23 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(227, 4).Value = $newDate

$ws.Cells.Item(227, 5).Value = 7
$ws.Cells.Item(227, 6).Value = 100114001
$ws.Cells.Item(227, 7).Value = 'Papa'
$ws.Cells.Item(227, 8).Value = 'Asterix'
$ws.Cells.Item(227, 9).Value = '1a nueva(o)'
$ws.Cells.Item(227, 10).Value = 1600
$ws.Cells.Item(227, 11).Value = 9000
$ws.Cells.Item(227, 12).Value = 9000
$ws.Cells.Item(227, 13).Value = 9000
$ws.Cells.Item(227, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(227, 15).Value = 'Región del Maule'
$ws.Cells.Item(227, 16).Value = 360
$ws.Cells.Item(227, 17).Value = 25
$ws.Cells.Item(227, 18).Value = 'Hortaliza'
